$wb = $excel.ActiveWorkbook

# --- times sheet: shift the scenario start/stop times ---
$wsTimes = $wb.Worksheets.Item("times")
$wsTimes.Range("B2").Value = 56249.99861111111
$wsTimes.Range("B3").Value = 56613.99861111111

# --- conventionals sheet: renumber + extend block list (AMIRIS 2.14 conventional fleet) ---
$wsConv = $wb.Worksheets.Item("conventionals")
$wsConv.Cells.Item(2, 1).Value = 0
$wsConv.Cells.Item(2, 2).Value = 99993300003
$wsConv.Cells.Item(2, 3).Value = "HYDROGEN"
$wsConv.Cells.Item(2, 4).Value = 1.5
$wsConv.Cells.Item(2, 5).Value = 0.4
$wsConv.Cells.Item(2, 6).Value = 1
$wsConv.Cells.Item(2, 7).Value = 1

$wsConv.Cells.Item(3, 1).Value = 1
$wsConv.Cells.Item(3, 2).Value = 99990600006
$wsConv.Cells.Item(3, 3).Value = "BIOMASS"
$wsConv.Cells.Item(3, 4).Value = 2.6
$wsConv.Cells.Item(3, 5).Value = 0.309
$wsConv.Cells.Item(3, 6).Value = 1
$wsConv.Cells.Item(3, 7).Value = 1

$wsConv.Cells.Item(4, 1).Value = 2
$wsConv.Cells.Item(4, 2).Value = 99991900008
$wsConv.Cells.Item(4, 3).Value = "NUCLEAR"
$wsConv.Cells.Item(4, 4).Value = 3.5
$wsConv.Cells.Item(4, 5).Value = 0.35
$wsConv.Cells.Item(4, 6).Value = 1
$wsConv.Cells.Item(4, 7).Value = 1

$wsConv.Cells.Item(5, 1).Value = 3
$wsConv.Cells.Item(5, 2).Value = 99992200009
$wsConv.Cells.Item(5, 3).Value = "NATURAL_GAS"
$wsConv.Cells.Item(5, 4).Value = 4.5
$wsConv.Cells.Item(5, 5).Value = 0.43
$wsConv.Cells.Item(5, 6).Value = 1
$wsConv.Cells.Item(5, 7).Value = 1

$wsConv.Cells.Item(6, 1).Value = 16
$wsConv.Cells.Item(6, 2).Value = 20253300023
$wsConv.Cells.Item(6, 3).Value = "HYDROGEN"
$wsConv.Cells.Item(6, 4).Value = 1.733432959545287
$wsConv.Cells.Item(6, 5).Value = 0.4
$wsConv.Cells.Item(6, 6).Value = 2000
$wsConv.Cells.Item(6, 7).Value = 2000

$wsConv.Cells.Item(7, 1).Value = 15
$wsConv.Cells.Item(7, 2).Value = 20273300024
$wsConv.Cells.Item(7, 3).Value = "HYDROGEN"
$wsConv.Cells.Item(7, 4).Value = 1.716227776089985
$wsConv.Cells.Item(7, 5).Value = 0.4
$wsConv.Cells.Item(7, 6).Value = 2000
$wsConv.Cells.Item(7, 7).Value = 2000

$wsConv.Cells.Item(8, 1).Value = 14
$wsConv.Cells.Item(8, 2).Value = 20293300025
$wsConv.Cells.Item(8, 3).Value = "HYDROGEN"
$wsConv.Cells.Item(8, 4).Value = 1.699193362629623
$wsConv.Cells.Item(8, 5).Value = 0.4
$wsConv.Cells.Item(8, 6).Value = 2000
$wsConv.Cells.Item(8, 7).Value = 2000

$wsConv.Cells.Item(9, 1).Value = 13
$wsConv.Cells.Item(9, 2).Value = 20313300026
$wsConv.Cells.Item(9, 3).Value = "HYDROGEN"
$wsConv.Cells.Item(9, 4).Value = 1.682328024187148
$wsConv.Cells.Item(9, 5).Value = 0.4
$wsConv.Cells.Item(9, 6).Value = 2000
$wsConv.Cells.Item(9, 7).Value = 2000

$wsConv.Cells.Item(10, 1).Value = 12
$wsConv.Cells.Item(10, 2).Value = 20333300027
$wsConv.Cells.Item(10, 3).Value = "HYDROGEN"
$wsConv.Cells.Item(10, 4).Value = 1.665630082608993
$wsConv.Cells.Item(10, 5).Value = 0.4
$wsConv.Cells.Item(10, 6).Value = 2000
$wsConv.Cells.Item(10, 7).Value = 2000

$wsConv.Cells.Item(11, 1).Value = 11
$wsConv.Cells.Item(11, 2).Value = 20353300028
$wsConv.Cells.Item(11, 3).Value = "HYDROGEN"
$wsConv.Cells.Item(11, 4).Value = 1.649097876398102
$wsConv.Cells.Item(11, 5).Value = 0.4
$wsConv.Cells.Item(11, 6).Value = 2000
$wsConv.Cells.Item(11, 7).Value = 2000

$wsConv.Cells.Item(12, 1).Value = 9
$wsConv.Cells.Item(12, 2).Value = 20393300030
$wsConv.Cells.Item(12, 3).Value = "HYDROGEN"
$wsConv.Cells.Item(12, 4).Value = 1.616524106382123
$wsConv.Cells.Item(12, 5).Value = 0.4
$wsConv.Cells.Item(12, 6).Value = 2000
$wsConv.Cells.Item(12, 7).Value = 2000

$wsConv.Cells.Item(13, 1).Value = 17
$wsConv.Cells.Item(13, 2).Value = 20233300022
$wsConv.Cells.Item(13, 3).Value = "HYDROGEN"
$wsConv.Cells.Item(13, 4).Value = 1.750810624964728
$wsConv.Cells.Item(13, 5).Value = 0.4
$wsConv.Cells.Item(13, 6).Value = 2000
$wsConv.Cells.Item(13, 7).Value = 2000

$wsConv.Cells.Item(14, 1).Value = 7
$wsConv.Cells.Item(14, 2).Value = 20413300031
$wsConv.Cells.Item(14, 3).Value = "HYDROGEN"
$wsConv.Cells.Item(14, 4).Value = 1.600479301385731
$wsConv.Cells.Item(14, 5).Value = 0.4
$wsConv.Cells.Item(14, 6).Value = 2000
$wsConv.Cells.Item(14, 7).Value = 2000

$wsConv.Cells.Item(15, 1).Value = 6
$wsConv.Cells.Item(15, 2).Value = 20433300032
$wsConv.Cells.Item(15, 3).Value = "HYDROGEN"
$wsConv.Cells.Item(15, 4).Value = 1.58459374905149
$wsConv.Cells.Item(15, 5).Value = 0.4
$wsConv.Cells.Item(15, 6).Value = 2000
$wsConv.Cells.Item(15, 7).Value = 2000

$wsConv.Cells.Item(16, 1).Value = 5
$wsConv.Cells.Item(16, 2).Value = 20453300033
$wsConv.Cells.Item(16, 3).Value = "HYDROGEN"
$wsConv.Cells.Item(16, 4).Value = 1.568865868717596
$wsConv.Cells.Item(16, 5).Value = 0.4
$wsConv.Cells.Item(16, 6).Value = 2000
$wsConv.Cells.Item(16, 7).Value = 2000

$wsConv.Cells.Item(17, 1).Value = 4
$wsConv.Cells.Item(17, 2).Value = 20473300034
$wsConv.Cells.Item(17, 3).Value = "HYDROGEN"
$wsConv.Cells.Item(17, 4).Value = 1.553294095411101
$wsConv.Cells.Item(17, 5).Value = 0.4
$wsConv.Cells.Item(17, 6).Value = 2000
$wsConv.Cells.Item(17, 7).Value = 2000

$wsConv.Cells.Item(18, 1).Value = 10
$wsConv.Cells.Item(18, 2).Value = 20373300029
$wsConv.Cells.Item(18, 3).Value = "HYDROGEN"
$wsConv.Cells.Item(18, 4).Value = 1.632729760548603
$wsConv.Cells.Item(18, 5).Value = 0.4
$wsConv.Cells.Item(18, 6).Value = 2000
$wsConv.Cells.Item(18, 7).Value = 2000

$wsConv.Cells.Item(19, 1).Value = 8
$wsConv.Cells.Item(19, 2).Value = 20401900035
$wsConv.Cells.Item(19, 3).Value = "NUCLEAR"
$wsConv.Cells.Item(19, 4).Value = 3.753123961749539
$wsConv.Cells.Item(19, 5).Value = 0.35
$wsConv.Cells.Item(19, 6).Value = 5000
$wsConv.Cells.Item(19, 7).Value = 5000

$wsConv.Cells.Item(20, 1).Value = 18
$wsConv.Cells.Item(20, 2).Value = 20151900036
$wsConv.Cells.Item(20, 3).Value = "NUCLEAR"
$wsConv.Cells.Item(20, 4).Value = 4.251522216620676
$wsConv.Cells.Item(20, 5).Value = 0.35
$wsConv.Cells.Item(20, 6).Value = 5000
$wsConv.Cells.Item(20, 7).Value = 5000

# --- renewables sheet: refreshed identifiers / LCOE figures ---
$wsRen = $wb.Worksheets.Item("renewables")
$wsRen.Cells.Item(2, 1).Value = 0
$wsRen.Cells.Item(2, 2).Value = 99990100002
$wsRen.Cells.Item(2, 3).Value = 1
$wsRen.Cells.Item(2, 4).Value = 0.5
$wsRen.Cells.Item(2, 5).Value = "WindOff"

$wsRen.Cells.Item(3, 1).Value = 1
$wsRen.Cells.Item(3, 2).Value = 99990300004
$wsRen.Cells.Item(3, 3).Value = 1
$wsRen.Cells.Item(3, 4).Value = 0.5
$wsRen.Cells.Item(3, 5).Value = "OtherPV"

$wsRen.Cells.Item(4, 1).Value = 2
$wsRen.Cells.Item(4, 2).Value = 99990200005
$wsRen.Cells.Item(4, 3).Value = 1
$wsRen.Cells.Item(4, 4).Value = 0.5
$wsRen.Cells.Item(4, 5).Value = "WindOn"

$wsRen.Cells.Item(5, 1).Value = 3
$wsRen.Cells.Item(5, 2).Value = 99990500007
$wsRen.Cells.Item(5, 3).Value = 1
$wsRen.Cells.Item(5, 4).Value = 0.5
$wsRen.Cells.Item(5, 5).Value = "PVRooftop"

$wsRen.Cells.Item(6, 1).Value = 4
$wsRen.Cells.Item(6, 2).Value = 20490200021
$wsRen.Cells.Item(6, 3).Value = 12000
$wsRen.Cells.Item(6, 4).Value = 0.5126256265640622
$wsRen.Cells.Item(6, 5).Value = "WindOn"

$wsRen.Cells.Item(7, 1).Value = 5
$wsRen.Cells.Item(7, 2).Value = 20490300037
$wsRen.Cells.Item(7, 3).Value = 35000
$wsRen.Cells.Item(7, 4).Value = 0.5126256265640622
$wsRen.Cells.Item(7, 5).Value = "OtherPV"

$wsRen.Cells.Item(8, 1).Value = 6
$wsRen.Cells.Item(8, 2).Value = 20400300038
$wsRen.Cells.Item(8, 3).Value = 35000
$wsRen.Cells.Item(8, 4).Value = 0.5361605659642198
$wsRen.Cells.Item(8, 5).Value = "OtherPV"

$wsRen.Cells.Item(9, 1).Value = 7
$wsRen.Cells.Item(9, 2).Value = 20370100040
$wsRen.Cells.Item(9, 3).Value = 25000
$wsRen.Cells.Item(9, 4).Value = 0.544243253516201
$wsRen.Cells.Item(9, 5).Value = "WindOff"

$wsRen.Cells.Item(10, 1).Value = 8
$wsRen.Cells.Item(10, 2).Value = 20250100039
$wsRen.Cells.Item(10, 3).Value = 25000
$wsRen.Cells.Item(10, 4).Value = 0.5778109865150957
$wsRen.Cells.Item(10, 5).Value = "WindOff"

# --- storages sheet: drop to 2 multi-agent storages, renumber ---
$wsStor = $wb.Worksheets.Item("storages")
$wsStor.Cells.Item(2, 1).Value = 0
$wsStor.Cells.Item(2, 2).Value = 99992600001
$wsStor.Cells.Item(2, 3).Value = "STORAGE"
$wsStor.Cells.Item(2, 4).Value = 4
$wsStor.Cells.Item(2, 5).Value = 0.92
$wsStor.Cells.Item(2, 6).Value = 0.92
$wsStor.Cells.Item(2, 7).Value = 0
$wsStor.Cells.Item(2, 8).Value = 1
$wsStor.Cells.Item(2, 9).Value = "MULTI_AGENT_MEDIAN"

$wsStor.Cells.Item(3, 1).Value = 1
$wsStor.Cells.Item(3, 2).Value = 20502600043
$wsStor.Cells.Item(3, 3).Value = "STORAGE"
$wsStor.Cells.Item(3, 4).Value = 4
$wsStor.Cells.Item(3, 5).Value = 0.92
$wsStor.Cells.Item(3, 6).Value = 0.92
$wsStor.Cells.Item(3, 7).Value = 0
$wsStor.Cells.Item(3, 8).Value = 2000
$wsStor.Cells.Item(3, 9).Value = "MULTI_AGENT_MEDIAN"

$wsStor.Range("A4:A10").EntireRow.Delete()

